$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (B:H) ---
$ws.Columns.Item(2).ColumnWidth = 54.166666666666664   # B -> 55
$ws.Columns.Item(3).ColumnWidth = 54.166666666666664   # C -> 55
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668   # D -> 30
$ws.Columns.Item(6).ColumnWidth = 16.166666666666668   # F -> 17
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666   # G -> 16
$ws.Columns.Item(8).ColumnWidth = 26.166666666666668   # H -> 27

# --- Helper to write a purely-numeric-looking ID as text (not a number) ---
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Row 2 ---
Set-TextValue $ws.Range("A2") "1328744"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328744"
$ws.Range("C2").Value = "[Remote] Front Desk – Customer Service (Semi Senior)"
$ws.Range("D2").Value = "No location available"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "18 applicants"
$ws.Range("G2").Value = "Remote"
$ws.Range("H2").Value = "Aurent LLC"

# --- Row 3 ---
Set-TextValue $ws.Range("A3") "1325846"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1325846"
$ws.Range("C3").Value = "Digital Media Strategist - Long Term"
$ws.Range("D3").Value = "Nugegoda, Sri Lanka"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "44 applicants"
$ws.Range("G3").Value = "6 - 18 Months"
$ws.Range("H3").Value = "Brand Corridor (Pvt) Ltd"

# --- Row 4 ---
Set-TextValue $ws.Range("A4") "1316788"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1316788"
$ws.Range("C4").Value = "Travel Coordinator"
$ws.Range("D4").Value = "Mexico City, CDMX, Mexico"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "128 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "Ikan Experience"

# --- Row 5 ---
Set-TextValue $ws.Range("A5") "1313206"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1313206"
$ws.Range("C5").Value = "Digital Media Strategist"
$ws.Range("D5").Value = "Colombo, Sri Lanka"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "47 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "Brand Corridor (Pvt) Ltd"

# --- Row 6 ---
Set-TextValue $ws.Range("A6") "1301868"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1301868"
$ws.Range("C6").Value = "Marketing Assistant"
$ws.Range("D6").Value = "Alor Setar, Kedah, Malaysia"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "435 applicants"
$ws.Range("G6").Value = "3 - 6 Months"
$ws.Range("H6").Value = "Yonhin Sdn. Bhd"

Write-Host "Edit complete"
